# Apply the cryptos.xlsx update (Tue Sep 24 08:14:10 UTC 2024, GitHub Actions run).
# Every changed cell holds plain text in the source workbook (t="inlineStr"), so we
# prefix values with a literal apostrophe to force Excel to store them as text even
# when they look numeric (e.g. "5.59"), then reset the cell style to "Normal" so the
# automatic quote-prefix formatting Excel applies is not left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

Set-TextValue "D2" "63.655.89"
Set-TextValue "E2" "  +0.35%  "
Set-TextValue "D3" "2.647.50"
Set-TextValue "E3" "  +0.23%  "
Set-TextValue "E4" "  +0.02%  "
Set-TextValue "D5" "602.39"
Set-TextValue "E5" "  +2.18%  "
Set-TextValue "D6" "147.37"
Set-TextValue "E6" "  +1.83%  "
Set-TextValue "E7" "  +0.03%  "
Set-TextValue "E8" "  +0.02%  "
Set-TextValue "E9" "  +1.50%  "
Set-TextValue "D10" "5.59"
Set-TextValue "E10" "  -0.88%  "
Set-TextValue "D11" "0.369"
Set-TextValue "E11" "  +4.53%  "
Set-TextValue "D13" "27.51"
Set-TextValue "E13" "  +0.06%  "
Set-TextValue "D14" "3.122.17"
Set-TextValue "E14" "  +0.24%  "
Set-TextValue "D15" "63.508.32"
Set-TextValue "E15" "  +0.17%  "
Set-TextValue "E16" "  +0.29%  "
Set-TextValue "D17" "2.647.86"
Set-TextValue "E17" "  +0.92%  "
Set-TextValue "D18" "11.50"
Set-TextValue "E18" "  +1.43%  "
Set-TextValue "E19" "  +4.27%  "
Set-TextValue "D20" "341.66"
Set-TextValue "E20" "  +0.49%  "
Set-TextValue "D21" "6.99"
Set-TextValue "E21" "  +4.52%  "
Set-TextValue "D23" "5.58"
Set-TextValue "E23" "  -3.11%  "
Set-TextValue "D24" "66.79"
Set-TextValue "E24" "  -1.14%  "
Set-TextValue "D25" "1.70"
Set-TextValue "E25" "  +0.10%  "
Set-TextValue "D26" "8.95"
Set-TextValue "E26" "  +6.19%  "
Set-TextValue "D27" "1.54"
Set-TextValue "E27" "  -0.64%  "
Set-TextValue "E28" "  -0.89%  "
Set-TextValue "D29" "547.50"
Set-TextValue "E29" "  -1.44%  "
Set-TextValue "D30" "0.999"
Set-TextValue "E30" "  -0.19%  "
Set-TextValue "D31" "7.84"
Set-TextValue "E31" "  +0.28%  "
Set-TextValue "D32" "2.06"
Set-TextValue "E32" "  +4.61%  "
Set-TextValue "E33" "  -2.85%  "
Set-TextValue "D34" "0.0₃0812"
Set-TextValue "E34" "  +1.28%  "
Set-TextValue "D35" "5.21"
Set-TextValue "E35" "  +7.09%  "
Set-TextValue "D36" "168.07"
Set-TextValue "E36" "  -3.96%  "
Set-TextValue "E37" "  +1.28%  "
Set-TextValue "B39" "Stacks"
Set-TextValue "C39" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D39" "1.91"
Set-TextValue "E39" "  +8.47%  "
Set-TextValue "B40" "EthereumClassic"
Set-TextValue "C40" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D40" "19.12"
Set-TextValue "E40" "  +0.47%  "
Set-TextValue "E41" "  -0.04%  "
Set-TextValue "D42" "169.07"
Set-TextValue "E42" "  -0.54%  "
Set-TextValue "B43" "Filecoin"
Set-TextValue "C43" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D43" "3.78"
Set-TextValue "E43" "  +1.97%  "
Set-TextValue "B44" "InjectiveProtocol"
Set-TextValue "C44" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D44" "22.72"
Set-TextValue "E44" "  +2.35%  "
Set-TextValue "D45" "0.0580"
Set-TextValue "E45" "  +5.73%  "
Set-TextValue "D46" "0.628"
Set-TextValue "E46" "  +0.08%  "
Set-TextValue "D47" "0.0248"
Set-TextValue "E47" "  +4.65%  "
Set-TextValue "E48" "  +0.62%  "
Set-TextValue "D49" "18.90"
Set-TextValue "E49" "  +1.17%  "
Set-TextValue "D50" "1.85"
Set-TextValue "E50" "  +9.14%  "
Set-TextValue "D51" "11.28"
Set-TextValue "E51" "  -0.65%  "
